$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to remain plain text
# (prevents Excel from auto-converting numeric-looking strings like
# "6.50" or "0.0234" into numbers and losing trailing/leading zeros,
# and preserves the cell's original style/formatting).
function Set-TextValue($ws, $addr, $value) {
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "60.636.06"
$ws.Range("E2").Value = "  -4.27%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "2.918.92"
$ws.Range("E3").Value = "  -3.66%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "529.69"
$ws.Range("E5").Value = "  -5.09%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "145.53"
$ws.Range("E6").Value = "  -6.28%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.11%  "

# Row 8 - XRP
Set-TextValue $ws "D8" "0.561"

# Row 9 - LidoStakedEther
Set-TextValue $ws "D9" "2.917.72"
$ws.Range("E9").Value = "  -4.03%  "

# Row 10 - Dogecoin
Set-TextValue $ws "D10" "0.109"
$ws.Range("E10").Value = "  -3.46%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -7.53%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -2.81%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D13" "3.429.40"
$ws.Range("E13").Value = "  -3.74%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.45%  "

# Row 15 - WrappedBTC
Set-TextValue $ws "D15" "60.753.99"
$ws.Range("E15").Value = "  -4.14%  "

# Row 16 - Avalanche
Set-TextValue $ws "D16" "23.04"
$ws.Range("E16").Value = "  -4.63%  "

# Row 17 - WrappedEther
Set-TextValue $ws "D17" "2.911.07"
$ws.Range("E17").Value = "  -4.04%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -6.13%  "

# Row 19 - Polkadot
Set-TextValue $ws "D19" "5.01"

# Row 20 - Chainlink
Set-TextValue $ws "D20" "11.69"
$ws.Range("E20").Value = "  -2.97%  "

# Row 21 - BitcoinCash
Set-TextValue $ws "D21" "365.43"
$ws.Range("E21").Value = "  -8.80%  "

# Row 22 - Uniswap
Set-TextValue $ws "D22" "6.50"
$ws.Range("E22").Value = "  -2.66%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.15%  "

# Row 24 - LEO
Set-TextValue $ws "D24" "5.68"
$ws.Range("E24").Value = "  -2.09%  "

# Row 25 - Litecoin
Set-TextValue $ws "D25" "64.72"
$ws.Range("E25").Value = "  -1.14%  "

# Row 26 - WrappedeETH
Set-TextValue $ws "D26" "3.054.54"
$ws.Range("E26").Value = "  -3.21%  "

# Row 27 - was Kaspa, now Polygon (rows 27/28 swapped)
$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws "D27" "0.455"
$ws.Range("E27").Value = "  -2.09%  "

# Row 28 - was Polygon, now Kaspa
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D28" "0.185"
$ws.Range("E28").Value = "  -2.90%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.30%  "

# Row 30 - PEPE
$sub3 = [string][char]0x2083
Set-TextValue $ws "D30" ("0.0" + $sub3 + "0878")
$ws.Range("E30").Value = "  -11.15%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue $ws "D31" "7.76"
$ws.Range("E31").Value = "  -11.05%  "

# Row 32 - USDe
$ws.Range("E32").Value = "  +0.05%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  -4.64%  "

# Row 34 - EthereumClassic
$ws.Range("E34").Value = "  -2.92%  "

# Row 35 - Monero
Set-TextValue $ws "D35" "158.10"
$ws.Range("E35").Value = "  -2.82%  "

# Row 36 - NEARProtocol
Set-TextValue $ws "D36" "4.43"
$ws.Range("E36").Value = "  -6.33%  "

# Row 37 - Aptos
Set-TextValue $ws "D37" "5.66"
$ws.Range("E37").Value = "  -6.50%  "

# Row 38 - Fetch.AI
$ws.Range("E38").Value = "  -9.48%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -6.86%  "

# Row 40 - OKB
Set-TextValue $ws "D40" "38.06"
$ws.Range("E40").Value = "  +0.57%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -6.16%  "

# Row 42 - Maker
Set-TextValue $ws "D42" "2.359.12"
$ws.Range("E42").Value = "  -7.28%  "

# Row 43 - Filecoin
Set-TextValue $ws "D43" "3.73"
$ws.Range("E43").Value = "  -5.71%  "

# Row 44 - Mantle
Set-TextValue $ws "D44" "0.648"
$ws.Range("E44").Value = "  -3.53%  "

# Row 45 - EnergySwap
Set-TextValue $ws "D45" "21.03"
$ws.Range("E45").Value = "  -8.46%  "

# Row 46 - Hedera
Set-TextValue $ws "D46" "0.0575"
$ws.Range("E46").Value = "  -4.23%  "

# Row 47 - FirstDigitalUSD
Set-TextValue $ws "D47" "0.998"
$ws.Range("E47").Value = "  +0.06%  "

# Row 48 - RenderToken
Set-TextValue $ws "D48" "4.96"
$ws.Range("E48").Value = "  -3.11%  "

# Row 49 - VeChain
Set-TextValue $ws "D49" "0.0234"
$ws.Range("E49").Value = "  -6.75%  "

# Row 50 - WhiteBITCoin
$ws.Range("E50").Value = "  -1.11%  "

# Row 51 - Stellar
Set-TextValue $ws "D51" "0.0930"
$ws.Range("E51").Value = "  -1.52%  "
